# "added spx500 for prop challenge"
#
# Adds a new "SPX500 (2017-2022)" block (rows 11-13) to Sheet1, mirroring the
# existing GBP block (rows 5-7): label in merged A11:A13, profit/bal dd/eq dd
# rows in B, raw C/D/E/F figures, and the two ratio helper formulas in H11 / J11.
# Also tidies up a couple of stray fill flags on K2/K8 and widens column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Write the two new ratio formulas FIRST, while their precedents (C11:C12,
#    E11:E12) are still untouched/General - this engine paints a formula
#    cell's number format from its precedents at assignment time, so doing
#    this before the precedent cells pick up the "#,##0" look keeps these two
#    cells on the plain General style the source rows actually use.
# ---------------------------------------------------------------------------
$ws.Range("H11").Formula = "=C11/C12"
$ws.Range("J11").Formula = "=E11/E12"

# Re-apply the correct look to those two cells (format-only paste, formulas
# untouched) now that the style has been pinned down.
$ws.Range("H5").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$ws.Range("J8").Copy() | Out-Null
$ws.Range("J11").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2) New data block, rows 11-13 (same visual style family as rows 5-7: A=GBP
#    label col, B/C/D/E/F = profit / bal dd / eq dd figures).
#    Merge the label column first so the merge doesn't go back and rewrite
#    the border edges of cells that already picked up their format.
# ---------------------------------------------------------------------------
$ws.Range("A11:A13").Merge() | Out-Null

# Row 11 ---------------------------------------------------------------
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = "SPX500 (2017-2022)"

$ws.Range("B5").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "profit"

$ws.Range("C5").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Value = 4451

$ws.Range("D5").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null         # stays blank

$ws.Range("E5").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Value = 69

$ws.Range("F5").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null         # stays blank

# Row 12 ---------------------------------------------------------------
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null         # blank, merged

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").Value = "bal dd"

$ws.Range("C6").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = 717

$ws.Range("D6").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null

$ws.Range("E6").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = 1.5

$ws.Range("F6").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null

# Row 13 ---------------------------------------------------------------
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null         # blank, merged

$ws.Range("B7").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = "eq dd"

$ws.Range("C7").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value = 912

$ws.Range("D7").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Range("E7").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Value = 307

$ws.Range("F7").Copy() | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null

# Two new (empty, just carrying a style) helper cells alongside the
# existing bal-dd / eq-dd rows of the block above.
$ws.Range("K2").Copy() | Out-Null
$ws.Range("J9").PasteSpecial(-4122) | Out-Null
$ws.Range("J10").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3) K2 / K8 had a stray "apply fill" flag with no actual fill behind it -
#    clear that so they match the plain centered style used elsewhere.
# ---------------------------------------------------------------------------
$ws.Range("K2").Interior.Pattern = -4142   # xlNone
$ws.Range("K8").Interior.Pattern = -4142   # xlNone

# ---------------------------------------------------------------------------
# 4) Column A needs to be wide enough for the new "SPX500 (2017-2022)" label.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 16.5

# ---------------------------------------------------------------------------
# 5) Leave the selection where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("B20").Select() | Out-Null
